# Added more margin to tabs.
#
# The flyer's "Tabs" group (inside the "Haiku" group on slide 1) holds ten
# vertical tab labels ("Tab 1" .. "Tab 10"). Give each of them a bit more
# top margin (0.2in -> 0.3in) so the text isn't crowded against the tab
# edge.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "Haiku" group shape that contains the tab textboxes.
$haiku = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Haiku") {
        $haiku = $shp
        break
    }
}

for ($i = 1; $i -le 10; $i++) {
    $tabName = "Tab $i"
    for ($j = 1; $j -le $haiku.GroupItems.Count; $j++) {
        $item = $haiku.GroupItems.Item($j)
        if ($item.Name -eq $tabName) {
            $item.TextFrame.MarginTop = 21.6
            break
        }
    }
}
